# Insert a new data row at row 268 (shifts existing rows 268..328 down to 269..329)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with the new record's data
$ws.Cells.Item(268, 1).Value = 6
$ws.Cells.Item(268, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(268, 3).Value = "Metropolitana"
$ws.Cells.Item(268, 4).Value = 44543
$ws.Cells.Item(268, 5).Value = 13
$ws.Cells.Item(268, 6).Value = 100112043
$ws.Cells.Item(268, 7).Value = "Pepino ensalada"
$ws.Cells.Item(268, 8).Value = "Sin especificar"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 2200
$ws.Cells.Item(268, 11).Value = 6000
$ws.Cells.Item(268, 12).Value = 7000
$ws.Cells.Item(268, 13).Value = 6455
$ws.Cells.Item(268, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(268, 15).Value = "Región Metropolitana"
$ws.Cells.Item(268, 16).Value = 129
$ws.Cells.Item(268, 17).Value = 50
$ws.Cells.Item(268, 18).Value = "Hortaliza"

# Match the date-column style (s="2", custom date number format) used by the other rows
$ws.Cells.Item(268, 4).NumberFormat = $ws.Cells.Item(269, 4).NumberFormat
